$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 110, pushing the existing rows 110-119 down to 111-120
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record
$ws.Range("A110").Value = 5
$ws.Range("B110").Value = "Macroferia Regional de Talca"
$ws.Range("C110").Value = "Maule"
$ws.Range("D110").Value = 45194
$ws.Range("E110").Value = 7
$ws.Range("F110").Value = 100112026
$ws.Range("G110").Value = "Haba"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 400
$ws.Range("K110").Value = 10000
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 10000
$ws.Range("N110").Value = "$/saco 25 kilos"
$ws.Range("O110").Value = "Región de O'Higgins"
$ws.Range("P110").Value = 400
$ws.Range("Q110").Value = 25
$ws.Range("R110").Value = "Hortaliza"

# Make sure the date cell keeps the same date number format as the other date cells in column D
$ws.Range("D110").NumberFormat = $ws.Range("D111").NumberFormat
